$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column E (header border/font) over to the two
# newly-added columns F and G so every cell picks up the same style (s=2)
# used across the rest of the header/data grid.
$ws.Range("E1:E17").Copy() | Out-Null
$ws.Range("F1:F17").PasteSpecial(-4122) | Out-Null
$ws.Range("E1:E17").Copy() | Out-Null
$ws.Range("G1:G17").PasteSpecial(-4122) | Out-Null

# New header labels for the added columns.
$ws.Range("F1").Value = "Nomor WA"
$ws.Range("G1").Value = "E-Mail"

# Widen the E:G block slightly (28 -> ~28.29 chars) to fit the new headers.
$ws.Range("E1:G17").ColumnWidth = 27.5

# Move/clear the active selection to the first data cell of the new block.
$ws.Range("F2").Select() | Out-Null
